$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final player/position/team table (post-edit), in display order.
$data = @(
    @("Dennis Schröder", "PG", "Golden State Warriors"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Terry Rozier", "PG", "Miami Heat"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Dereck Lively II", "C", "Dallas Mavericks"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors")
)

# Clear out the old data range (header stays at row 1; old table ran to row 18).
$ws.Range("A2:C18").ClearContents()

# Write the new table starting at row 2.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
